$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Catchment_LanduseMix")

# Update fD values in column D (D2:D4) to be area-proportional (all set to 0)
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0

# Update C2 from 40 to 30
$ws.Range("C2").Value = 30

# Add new landuse type "main_road" as row 5
$ws.Range("A5").Value = "main_road"
$ws.Range("B5").Value = 0.6
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 0

# Copy style from row 4 to row 5 so formatting matches (fill/border)
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C4:D4").Copy()
$ws.Range("C5:D5").PasteSpecial(-4122)  # xlPasteFormats

# Update the selection to match B6 (as shown in diff)
$ws.Range("B6").Select()
